# 22 march changes, fixed plots
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix up existing rows ---

# Row 7: project_subtype mislabeled as "biostatistics"/"Data processing" -> "Biostatistics"
$ws.Range("B7").Value = "Biostatistics"
$ws.Range("C7").Value = "Biostatistics"

# Row 8/9: clarify remarks wording
$ws.Range("I8").Value = "There were 2 versions of this analysis (i.e. different models)"
$ws.Range("I9").Value = "There were 4 versions of this analysis (i.e. different models)"

# Row 14: project_type/project_subtype capitalization fix
$ws.Range("B14").Value = "Biostatistics"
$ws.Range("C14").Value = "Biostatistics"

# Row 20: project_type corrected to Bulk RNAseq
$ws.Range("B20").Value = "Bulk RNAseq"

# Row 22: project_subtype corrected to Benchmarking
$ws.Range("C22").Value = "Benchmarking"

# Row 25: add missing remarks
$ws.Range("I25").Value = "Benchmarking in 1 dataset"

# Row 26: project_type corrected to WES CN
$ws.Range("B26").Value = "WES CN"

# Row 27: project_subtype + remarks
$ws.Range("C27").Value = "Regulatory Analysis"
$ws.Range("I27").Value = "Calra requested help in interpreting gene regulatory information with respect to her gene of interest, TRAIL"

# Row 28: add missing remarks
$ws.Range("I28").Value = "Qun is benchmarking Single Cell Data and needed help with a script to extract the necessary files"

# Row 29: completed flag + remarks
$ws.Range("H29").Value = "Y"
$ws.Range("I29").Value = "Benchmarking in 1 dataset"

# Row 30: add missing remarks
$ws.Range("I30").Value = "Benchmarking PERK in T-cells in two different datasets"

# --- New row 31 ---
$ws.Range("A31").Value = "Benchmarking Publically Available Methylation Datasets"
$ws.Range("B31").Value = "Epigenetics"
$ws.Range("C31").Value = "Benchmarking"
$ws.Range("D31").Value = "Ewout Landeloos"
$ws.Range("E31").Value = "Marine"

# Copy the date format (numFmtId 14, m/d/yyyy) from an existing date cell
# so the new cell reuses the same style instead of minting a new number format.
$ws.Range("F25").Copy()
$ws.Range("F31").PasteSpecial(-4122)
$newDate = Get-Date -Year 2021 -Month 3 -Day 22 -Hour 0 -Minute 0 -Second 0
$ws.Range("F31").Value = $newDate

$ws.Range("G31").Value = 1
$ws.Range("H31").Value = "Unfinished"
$ws.Range("I31").Value = "Ewout wants to test different methylation datasets and compare the results with his data"

# Update selection to match the saved workbook state
$ws.Range("C32").Select()
